# Fix typos in sample arc: rename the "C" (collision-prone) sample series
# to "CC" so it no longer clashes with the "Co" series, and correct the
# row-4 typo so it reads "CC3_prep" instead of a duplicated "C2_prep".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MS")

$ws.Range("A2").Value = "CC1_prep"
$ws.Range("A3").Value = "CC2_prep"
$ws.Range("A4").Value = "CC3_prep"
$ws.Range("A5").Value = "Co1_prep"
$ws.Range("A6").Value = "Co2_prep"
$ws.Range("A7").Value = "Co3_prep"

$ws.Range("AI2").Value = "CC1_measured"
$ws.Range("AI3").Value = "CC2_measured"
$ws.Range("AI4").Value = "CC2_measured"
$ws.Range("AI5").Value = "Co1_measured"
$ws.Range("AI6").Value = "Co2_measured"
$ws.Range("AI7").Value = "Co3_measured"

# Leave the selection where the author's last edit was (AI4, the row-4
# typo fix), scrolled right so column Q is the left-most visible column.
$ws.Activate() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollColumn = 17
$ws.Range("AI4").Select() | Out-Null
